$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.073442253343571906,
    -0.0099999993202786186,
    -0.0089999993301272951,
    0.28398601043740612,
    -0.0059999993539161522,
    -0.0059999993321717682,
    -0.019999999202724439,
    -0.019999999199528773,
    -0.0059999993256463213,
    -0.0059999993247572547,
    0.053792536330671936,
    -0.005999999322115368,
    -0.0065524171825410704,
    -0.01199999925252726,
    -0.0059999993048709399,
    -0.0059999993024986154,
    -0.0059999992993624574,
    -0.0089999992711806698,
    -0.008999999327577779,
    -0.0089999993090827957,
    -0.054310343679125239,
    -0.0089999993045828397,
    -0.0089999993241285381,
    -0.041999999011197886,
    -0.041999999005809308,
    -0.0059999993294930221,
    -0.0059999993259536311,
    -0.0059999993117347827,
    -0.011999999246096849,
    -0.019999999167885196,
    -0.014999999206350978,
    -0.02099999915024231,
    -0.0059999992871206942
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
